$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '309.24'

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '0.39%'

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '41.11'

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '0.36%'

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '5.217'

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '2.26%'

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '0.07691'

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '0.88%'

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '1.644'

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '2.34%'

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.9146'

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '1.36%'

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.1246'

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '11.14%'

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '2.48%'

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.09200'

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '0.75%'

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.04215'

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '-0.27%'

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.1051'

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '0.001261'

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '0.01%'

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.005753'

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '-0.16%'

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '1,903.19%'

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '-0.04%'

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '1.58%'

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '7.404'

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '11.77%'

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '1.18%'

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '0.58%'

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '0.04035'

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '-0.83%'

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '0.001265'

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '1.71%'

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.004092'

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '-0.51%'

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '0.0001302'

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '0.01%'

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.02569'

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '7.79%'

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.05338'

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '3.07%'

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.007839'

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '0.71%'

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '1.24%'

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.006664'

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '-5.59%'

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.001862'

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '-4.60%'

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.008050'

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '4.10%'

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.3070'

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '-0.33%'

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.00006729'

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '-3.97%'

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.00000000751'

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '0.01%'

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.2987'

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '436.50%'

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.00002102'

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '0.01%'

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.0002002'

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '0.01%'
